# Fix the "Faculty Cumulative Metric Score" table: correct the Wi
# weightage values for the Academic / Professional Activities columns
# and the Score / Weighted Score figures that are derived from them.

$d = $word.ActiveDocument

# The table is the 3rd table in the document:
#   Row 1: title "Faculty Cumulative Metric Score"
#   Row 2: column headers (1.Academic, 2.Professional Activities, ...)
#   Row 3: Score(S)
#   Row 4: Wi
#   Row 5: Weighted Score
$table = $d.Tables.Item(3)

function Set-CellText($row, $col, $newText) {
    $cellRange = $table.Cell($row, $col).Range
    # A table cell's Range includes the trailing paragraph mark (Chr 13)
    # and cell mark (Chr 7); drop the cell mark so we only overwrite the
    # actual cell content and keep the cell's own paragraph intact.
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $newText
}

# Row 3 (Score(S)) - Total Marks column, recalculated
Set-CellText 3 7 "11.35"

# Row 4 (Wi) - weightage factors
Set-CellText 4 2 "0.2"
Set-CellText 4 3 "0.4"

# Row 5 (Weighted Score) - recalculated from the corrected Wi values
Set-CellText 5 2 "3.6"
Set-CellText 5 3 "6.4"
